$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 65167
$ws.Cells.Item(2, 5).Value = 1285856360537
$ws.Cells.Item(2, 6).Value = 21650378234
$ws.Cells.Item(2, 7).Value = 2.38921

# Row 3
$ws.Cells.Item(3, 4).Value = 3162.43
$ws.Cells.Item(3, 5).Value = 380656097535
$ws.Cells.Item(3, 6).Value = 10460133831
$ws.Cells.Item(3, 7).Value = 3.81157

# Row 4
$ws.Cells.Item(4, 4).Value = 0.999634
$ws.Cells.Item(4, 5).Value = 109837986989
$ws.Cells.Item(4, 6).Value = 35814660588
$ws.Cells.Item(4, 7).Value = -0.14104

# Row 5
$ws.Cells.Item(5, 4).Value = 577.04
$ws.Cells.Item(5, 5).Value = 88997377881
$ws.Cells.Item(5, 6).Value = 1023560824
$ws.Cells.Item(5, 7).Value = 3.73376

# Row 6
$ws.Cells.Item(6, 4).Value = 150.38
$ws.Cells.Item(6, 5).Value = 67323891460
$ws.Cells.Item(6, 6).Value = 3154837312
$ws.Cells.Item(6, 7).Value = 6.31147

# Row 7
$ws.Cells.Item(7, 4).Value = 0.999115
$ws.Cells.Item(7, 5).Value = 33921182729
$ws.Cells.Item(7, 6).Value = 4967071302
$ws.Cells.Item(7, 7).Value = -0.05891

# Row 8
$ws.Cells.Item(8, 4).Value = 3161.21
$ws.Cells.Item(8, 5).Value = 29573964223
$ws.Cells.Item(8, 6).Value = 83475183
$ws.Cells.Item(8, 7).Value = 3.84677

# Row 9
$ws.Cells.Item(9, 4).Value = 0.528042
$ws.Cells.Item(9, 5).Value = 29179386461
$ws.Cells.Item(9, 6).Value = 919870439
$ws.Cells.Item(9, 7).Value = 2.02663

# Row 10
$ws.Cells.Item(10, 2).Value = "DOGE"
$ws.Cells.Item(10, 3).Value = "Dogecoin"
$ws.Cells.Item(10, 4).Value = 0.161294
$ws.Cells.Item(10, 5).Value = 23338959876
$ws.Cells.Item(10, 6).Value = 1656714530
$ws.Cells.Item(10, 7).Value = 6.07605

# Row 11
$ws.Cells.Item(11, 2).Value = "TON"
$ws.Cells.Item(11, 3).Value = "Toncoin"
$ws.Cells.Item(11, 4).Value = 6.16
$ws.Cells.Item(11, 5).Value = 21433855756
$ws.Cells.Item(11, 6).Value = 263750603
$ws.Cells.Item(11, 7).Value = 0.25616

# Row 12
$ws.Cells.Item(12, 4).Value = 0.5004150000000001
$ws.Cells.Item(12, 5).Value = 17732845591
$ws.Cells.Item(12, 6).Value = 404218775
$ws.Cells.Item(12, 7).Value = 4.75736

# Row 13
$ws.Cells.Item(13, 2).Value = "SHIB"
$ws.Cells.Item(13, 3).Value = "Shiba Inu"
$ws.Cells.Item(13, 4).Value = 0.0000268
$ws.Cells.Item(13, 5).Value = 15808645011
$ws.Cells.Item(13, 6).Value = 1152002688
$ws.Cells.Item(13, 7).Value = 16.35217

# Row 14
$ws.Cells.Item(14, 2).Value = "AVAX"
$ws.Cells.Item(14, 3).Value = "Avalanche"
$ws.Cells.Item(14, 4).Value = 37.31
$ws.Cells.Item(14, 5).Value = 14152141005
$ws.Cells.Item(14, 6).Value = 487064009
$ws.Cells.Item(14, 7).Value = 6.46406

# Row 15
$ws.Cells.Item(15, 4).Value = 65236
$ws.Cells.Item(15, 5).Value = 10131867455
$ws.Cells.Item(15, 6).Value = 163390002
$ws.Cells.Item(15, 7).Value = 2.32286

# Row 16
$ws.Cells.Item(16, 2).Value = "BCH"
$ws.Cells.Item(16, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(16, 4).Value = 510.83
$ws.Cells.Item(16, 5).Value = 10097849184
$ws.Cells.Item(16, 6).Value = 368311916
$ws.Cells.Item(16, 7).Value = 5.21581

# Row 17
$ws.Cells.Item(17, 2).Value = "TRX"
$ws.Cells.Item(17, 3).Value = "TRON"
$ws.Cells.Item(17, 4).Value = 0.110969
$ws.Cells.Item(17, 5).Value = 9722501396
$ws.Cells.Item(17, 6).Value = 278348564
$ws.Cells.Item(17, 7).Value = 1.28244

# Row 18
$ws.Cells.Item(18, 4).Value = 7.15
$ws.Cells.Item(18, 5).Value = 9704039943
$ws.Cells.Item(18, 6).Value = 185315831
$ws.Cells.Item(18, 7).Value = 6.10454

# Row 19
$ws.Cells.Item(19, 4).Value = 14.84
$ws.Cells.Item(19, 5).Value = 8742538777
$ws.Cells.Item(19, 6).Value = 311903752
$ws.Cells.Item(19, 7).Value = 4.81562

# Row 20
$ws.Cells.Item(20, 2).Value = "ICP"
$ws.Cells.Item(20, 3).Value = "Internet Computer"
$ws.Cells.Item(20, 4).Value = 15.25
$ws.Cells.Item(20, 5).Value = 7088337554
$ws.Cells.Item(20, 6).Value = 230990072
$ws.Cells.Item(20, 7).Value = 5.6414

# Row 21
$ws.Cells.Item(21, 2).Value = "NEAR"
$ws.Cells.Item(21, 3).Value = "NEAR Protocol"
$ws.Cells.Item(21, 4).Value = 6.31
$ws.Cells.Item(21, 5).Value = 6737443843
$ws.Cells.Item(21, 6).Value = 392598059
$ws.Cells.Item(21, 7).Value = 11.21983

# Row 22
$ws.Cells.Item(22, 2).Value = "MATIC"
$ws.Cells.Item(22, 3).Value = "Polygon"
$ws.Cells.Item(22, 4).Value = 0.722213
$ws.Cells.Item(22, 5).Value = 6726904459
$ws.Cells.Item(22, 6).Value = 294092710
$ws.Cells.Item(22, 7).Value = 5.81374

# Row 23
$ws.Cells.Item(23, 2).Value = "LTC"
$ws.Cells.Item(23, 3).Value = "Litecoin"
$ws.Cells.Item(23, 4).Value = 84.92
$ws.Cells.Item(23, 5).Value = 6338948133
$ws.Cells.Item(23, 6).Value = 347638758
$ws.Cells.Item(23, 7).Value = 3.0582

# Row 24
$ws.Cells.Item(24, 2).Value = "UNI"
$ws.Cells.Item(24, 3).Value = "Uniswap"
$ws.Cells.Item(24, 4).Value = 7.79
$ws.Cells.Item(24, 5).Value = 5878470129
$ws.Cells.Item(24, 6).Value = 125141781
$ws.Cells.Item(24, 7).Value = 3.8001

# Row 25
$ws.Cells.Item(25, 2).Value = "LEO"
$ws.Cells.Item(25, 3).Value = "LEO Token"
$ws.Cells.Item(25, 4).Value = 5.75
$ws.Cells.Item(25, 5).Value = 5325304943
$ws.Cells.Item(25, 6).Value = 1123640
$ws.Cells.Item(25, 7).Value = 0.20603

# Row 26
$ws.Cells.Item(26, 4).Value = 1.002
$ws.Cells.Item(26, 5).Value = 5148951014
$ws.Cells.Item(26, 6).Value = 496740453
$ws.Cells.Item(26, 7).Value = 0.14206

# Row 27
$ws.Cells.Item(27, 2).Value = "APT"
$ws.Cells.Item(27, 3).Value = "Aptos"
$ws.Cells.Item(27, 4).Value = 9.98
$ws.Cells.Item(27, 5).Value = 4269268905
$ws.Cells.Item(27, 6).Value = 123275240
$ws.Cells.Item(27, 7).Value = 4.167

# Row 28
$ws.Cells.Item(28, 2).Value = "STX"
$ws.Cells.Item(28, 3).Value = "Stacks"
$ws.Cells.Item(28, 4).Value = 2.83
$ws.Cells.Item(28, 5).Value = 4121088813
$ws.Cells.Item(28, 6).Value = 103656643
$ws.Cells.Item(28, 7).Value = 15.16295

# Row 29
$ws.Cells.Item(29, 2).Value = "ETC"
$ws.Cells.Item(29, 3).Value = "Ethereum Classic"
$ws.Cells.Item(29, 4).Value = 27.75
$ws.Cells.Item(29, 5).Value = 4089016384
$ws.Cells.Item(29, 6).Value = 174120125
$ws.Cells.Item(29, 7).Value = 5.21342

# Row 30
$ws.Cells.Item(30, 2).Value = "MNT"
$ws.Cells.Item(30, 3).Value = "Mantle"
$ws.Cells.Item(30, 4).Value = 1.2
$ws.Cells.Item(30, 5).Value = 3919913297
$ws.Cells.Item(30, 6).Value = 49545631
$ws.Cells.Item(30, 7).Value = 4.27778

# Row 31
$ws.Cells.Item(31, 2).Value = "FDUSD"
$ws.Cells.Item(31, 3).Value = "First Digital USD"
$ws.Cells.Item(31, 4).Value = 0.996842
$ws.Cells.Item(31, 5).Value = 3590329941
$ws.Cells.Item(31, 6).Value = 5563716800
$ws.Cells.Item(31, 7).Value = -0.27551

# Row 32
$ws.Cells.Item(32, 2).Value = "FIL"
$ws.Cells.Item(32, 3).Value = "Filecoin"
$ws.Cells.Item(32, 4).Value = 6.55
$ws.Cells.Item(32, 5).Value = 3566125827
$ws.Cells.Item(32, 6).Value = 199944914
$ws.Cells.Item(32, 7).Value = 6.07506

# Row 33
$ws.Cells.Item(33, 2).Value = "RNDR"
$ws.Cells.Item(33, 3).Value = "Render"
$ws.Cells.Item(33, 4).Value = 9.039999999999999
$ws.Cells.Item(33, 5).Value = 3486848953
$ws.Cells.Item(33, 6).Value = 169894511
$ws.Cells.Item(33, 7).Value = 11.697

# Row 34
$ws.Cells.Item(34, 2).Value = "CRO"
$ws.Cells.Item(34, 3).Value = "Cronos"
$ws.Cells.Item(34, 4).Value = 0.130083
$ws.Cells.Item(34, 5).Value = 3482562232
$ws.Cells.Item(34, 6).Value = 16021315
$ws.Cells.Item(34, 7).Value = 5.22879

# Row 35
$ws.Cells.Item(35, 2).Value = "ATOM"
$ws.Cells.Item(35, 3).Value = "Cosmos Hub"
$ws.Cells.Item(35, 4).Value = 8.66
$ws.Cells.Item(35, 5).Value = 3393373696
$ws.Cells.Item(35, 6).Value = 112957758
$ws.Cells.Item(35, 7).Value = 4.11151

# Row 36
$ws.Cells.Item(36, 2).Value = "XLM"
$ws.Cells.Item(36, 3).Value = "Stellar"
$ws.Cells.Item(36, 4).Value = 0.114555
$ws.Cells.Item(36, 5).Value = 3317003152
$ws.Cells.Item(36, 6).Value = 56946743
$ws.Cells.Item(36, 7).Value = 0.95014

# Row 37
$ws.Cells.Item(37, 4).Value = 55.19
$ws.Cells.Item(37, 5).Value = 3315645562
$ws.Cells.Item(37, 6).Value = 10160155
$ws.Cells.Item(37, 7).Value = -0.332

# Row 38
$ws.Cells.Item(38, 2).Value = "HBAR"
$ws.Cells.Item(38, 3).Value = "Hedera"
$ws.Cells.Item(38, 4).Value = 0.090183
$ws.Cells.Item(38, 5).Value = 3245321661
$ws.Cells.Item(38, 6).Value = 47450271
$ws.Cells.Item(38, 7).Value = 11.16918

# Row 39
$ws.Cells.Item(39, 2).Value = "ARB"
$ws.Cells.Item(39, 3).Value = "Arbitrum"
$ws.Cells.Item(39, 4).Value = 1.2
$ws.Cells.Item(39, 5).Value = 3199726355
$ws.Cells.Item(39, 6).Value = 212985096
$ws.Cells.Item(39, 7).Value = 6.11965

# Row 40
$ws.Cells.Item(40, 2).Value = "IMX"
$ws.Cells.Item(40, 3).Value = "Immutable"
$ws.Cells.Item(40, 4).Value = 2.18
$ws.Cells.Item(40, 5).Value = 3185371557
$ws.Cells.Item(40, 6).Value = 53228825
$ws.Cells.Item(40, 7).Value = 7.53123

# Row 41
$ws.Cells.Item(41, 2).Value = "EZETH"
$ws.Cells.Item(41, 3).Value = "Renzo Restaked ETH"
$ws.Cells.Item(41, 4).Value = 3189.96
$ws.Cells.Item(41, 5).Value = 3183532399
$ws.Cells.Item(41, 6).Value = 79980860
$ws.Cells.Item(41, 7).Value = 3.62206

# Row 42
$ws.Cells.Item(42, 2).Value = "TAO"
$ws.Cells.Item(42, 3).Value = "Bittensor"
$ws.Cells.Item(42, 4).Value = 473.03
$ws.Cells.Item(42, 5).Value = 3151141869
$ws.Cells.Item(42, 6).Value = 32710880
$ws.Cells.Item(42, 7).Value = 7.7267

# Row 43
$ws.Cells.Item(43, 2).Value = "VET"
$ws.Cells.Item(43, 3).Value = "VeChain"
$ws.Cells.Item(43, 4).Value = 0.04224241
$ws.Cells.Item(43, 5).Value = 3080340150
$ws.Cells.Item(43, 6).Value = 69364566
$ws.Cells.Item(43, 7).Value = 3.454

# Row 44
$ws.Cells.Item(44, 4).Value = 3.04
$ws.Cells.Item(44, 5).Value = 3050200015
$ws.Cells.Item(44, 6).Value = 601623097
$ws.Cells.Item(44, 7).Value = 10.44023

# Row 45
$ws.Cells.Item(45, 4).Value = 3058.8
$ws.Cells.Item(45, 5).Value = 2837429305
$ws.Cells.Item(45, 6).Value = 104402686
$ws.Cells.Item(45, 7).Value = 0.49827

# Row 46
$ws.Cells.Item(46, 2).Value = "KAS"
$ws.Cells.Item(46, 3).Value = "Kaspa"
$ws.Cells.Item(46, 4).Value = 0.117434
$ws.Cells.Item(46, 5).Value = 2716928958
$ws.Cells.Item(46, 6).Value = 33345159
$ws.Cells.Item(46, 7).Value = 1.74528

# Row 47
$ws.Cells.Item(47, 4).Value = 0.28424
$ws.Cells.Item(47, 5).Value = 2700089429
$ws.Cells.Item(47, 6).Value = 95829464
$ws.Cells.Item(47, 7).Value = 4.70704

# Row 48
$ws.Cells.Item(48, 2).Value = "INJ"
$ws.Cells.Item(48, 3).Value = "Injective"
$ws.Cells.Item(48, 4).Value = 28.99
$ws.Cells.Item(48, 5).Value = 2619707541
$ws.Cells.Item(48, 6).Value = 121917233
$ws.Cells.Item(48, 7).Value = 5.06998

# Row 49
$ws.Cells.Item(49, 2).Value = "OP"
$ws.Cells.Item(49, 3).Value = "Optimism"
$ws.Cells.Item(49, 4).Value = 2.47
$ws.Cells.Item(49, 5).Value = 2589251372
$ws.Cells.Item(49, 6).Value = 206107489
$ws.Cells.Item(49, 7).Value = 8.021940000000001

# Row 50
$ws.Cells.Item(50, 2).Value = "FET"
$ws.Cells.Item(50, 3).Value = "Fetch.ai"
$ws.Cells.Item(50, 4).Value = 2.45
$ws.Cells.Item(50, 5).Value = 2573943317
$ws.Cells.Item(50, 6).Value = 369802036
$ws.Cells.Item(50, 7).Value = 10.52097

# Row 51
$ws.Cells.Item(51, 2).Value = "PEPE"
$ws.Cells.Item(51, 3).Value = "Pepe"
$ws.Cells.Item(51, 4).Value = 0.00000601
$ws.Cells.Item(51, 5).Value = 2532730982
$ws.Cells.Item(51, 6).Value = 723356028
$ws.Cells.Item(51, 7).Value = 17.65564
